$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new student row (id 2) to the register
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Marcus Thomas"
$ws.Range("C3").Value = "Thomas"
$ws.Range("D3").Value = "Marcus"
$ws.Range("E3").Value = "Paris HEC"

# Reflect the last active selection after entering the new row
$ws.Range("E4").Select()
